# "updated feasible set and added inputs"
# Lower the feasible "Max" charge bound for Lithium, Sodium and Potassium
# from 2 to 1, then leave the active selection on D6 (next empty Max cell)
# to invite further input.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = 1
$ws.Range("D4").Value = 1
$ws.Range("D5").Value = 1

$ws.Range("D6").Select()
